$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with refreshed market data ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.136.94'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.622.37'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.28'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0846'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.621.48'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.085.52'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.58'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0742'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.87'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.24%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.96%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.48%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.57'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0506'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.99'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.339.52'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.95%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.859'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.802'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '65.53'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.58%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.23'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.758.99'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.913'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +36.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.38'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.42%  '

# --- Rows 48-51: list refreshed, BabyDogeCoin dropped off and USDD entered; ranks shifted up ---
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0513'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.37%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0995'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.78%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.56'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.08%  '
